# Doing Updates for Financials
# A new "D" column (latest fiscal-year period) is inserted into the SJT sheet,
# pushing the existing D:K data one column to the right (E:L) and populating
# the new column D with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before column D. Excel automatically shifts the
#    existing D:K cell contents (values, shared-string refs, styles) to E:L
#    and grows the sheet dimension from K102 to L102.
$ws.Columns("D").Insert()

# 2) The newly inserted column D cells have no explicit style yet (they pick
#    up the worksheet/column default). Copy the number formats/fonts from the
#    (now shifted) column E, which still carries the correct per-row style
#    (date format for header rows, right-aligned number format for data rows).
#    Only the row ranges that actually hold table data are touched, so blank
#    label-only rows (5, 6, 37, 79) do not gain a stray column-D cell.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3) Populate the new column D with the latest period's values.

# 3a. Period-ending dates (header rows for each of the three statements)
$ws.Range("D7").Value  = 43465
$ws.Range("D38").Value = 43465
$ws.Range("D80").Value = 43465

# 3b. Rows whose new-period figure is "NA" (reuses existing shared string)
$naRows = 9,10,12,14
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 4).Value = "NA"
}

# 3c. Rows whose new-period figure is 0
$zeroRows = 13,15,20,22,24,25,28,29,30,31,32,34,42,43,44,45,46,47,49,50,51,52,53,57,58,60,61,62,63,64,65,68,69,70,71,72,73,74,75,77,83,84,85,86,87,88,91,92,93,94,97,98,99,101
foreach ($r in $zeroRows) {
    $ws.Cells.Item($r, 4).Value = 0
}

# 3d. Rows that stay blank in the new period (section spacer / subtotal rows
#     with no values at all) - nothing to do, they already inherited the
#     correct blank style from the format copy above.

# 3e. Rows with specific new numeric figures for the latest period
$valueMap = @{
    8   = 19500
    17  = 1500
    18  = 18000
    21  = 18000
    23  = 18000
    26  = 18000
    27  = 18000
    33  = 18000
    35  = 18000
    41  = 2100
    48  = 5800
    54  = 8000
    59  = 1100
    66  = 2100
    76  = 5800
    81  = 18000
    89  = 17200
    96  = -18000
    100 = -18000
    102 = -700
}
foreach ($r in $valueMap.Keys) {
    $ws.Cells.Item($r, 4).Value = $valueMap[$r]
}
